# Add a "Turkey" worksheet (Zettler Turkey test data) based on the
# existing "Spain" sheet, matching the layout used by the other
# country sheets in this workbook.

$wb = $excel.ActiveWorkbook

$spain = $wb.Worksheets.Item("Spain")

# Duplicate the Spain sheet right after itself; the copy becomes the
# new active sheet, placed as the last tab in the workbook.
$spain.Copy($null, $spain)
$turkey = $wb.Worksheets.Item($wb.Worksheets.Count)
$turkey.Name = "Turkey"

# Fill in the Turkey-specific market name and Jira ticket reference.
$turkey.Range("B2").Value = "Turkey Market"
$turkey.Range("B4").Value = "NGC-3191/T3310"

# The copied rows 3-5 inherited Spain's taller row height (needed for
# its longer wrapped text); auto-fit them back down to the standard
# height used on the other country sheets.
$turkey.Rows.Item(3).AutoFit() | Out-Null
$turkey.Rows.Item(4).AutoFit() | Out-Null
$turkey.Rows.Item(5).AutoFit() | Out-Null

# Column C keeps the default sheet width instead of Spain's custom one,
# and column D is widened a bit for the new sheet.
$turkey.Columns.Item(3).ColumnWidth = $turkey.StandardWidth
$turkey.Columns.Item(4).ColumnWidth = 30

# Leave a fresh selection on the new sheet.
$turkey.Range("E15").Select()

# Clear the old selection/active-cell state left on Spain and select
# its full data range instead.
$spain.Range("A1:D10").Select()

# Make Turkey (the newly added, last tab) the active sheet/tab.
$turkey.Activate()
